$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2-19 changes from 2023-09-16 (45185) to 2023-10-05 (45204)
for ($row = 2; $row -le 19; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
